$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row 68 with the corrected lat/lon for "Near Mumbai coast"
$ws.Range("A68").Value = "Near Mumbai coast"
$ws.Range("B68").Value = 19.0847934780103
$ws.Range("C68").Value = 72.870532634994

# Extend the D column formula (Port~lat~long) down into the new row,
# the same way D67 extends the D3:D66 shared formula pattern
$ws.Range("D67:D68").FormulaR1C1 = "=RC[-3]&""~""&RC[-2]&""~""&RC[-1]"

# Match the scrolled/selected view state left behind by the edit
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 55
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A74").Select() | Out-Null
